$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("B6").Value = 'F001'
$ws.Range("C6").Value = 'Application is installed and launched'
$ws.Range("D6").Value = 'Network Packet Capture - Start'
$ws.Range("E6").Value = '1. Navigate to Network Packet Capture section<br>2. Click on Start button'
$ws.Range("F6").Value = 'Packet capture starts successfully'
$ws.Range("G6").Value = ''
$ws.Range("H6").Value = ''

# Row 7
$ws.Range("B7").Value = 'F002'
$ws.Range("C7").Value = 'Packet capture is in progress'
$ws.Range("D7").Value = 'Network Packet Capture - Stop'
$ws.Range("E7").Value = '1. Navigate to Network Packet Capture section<br>2. Click on Stop button'
$ws.Range("F7").Value = '1. Packet capture stops<br>2. .pcap file is generated<br>3. File is copied to MFP''s Shared Folder<br>4. Shared Folder opens automatically'
$ws.Range("G7").Value = ''
$ws.Range("H7").Value = ''

# Row 8
$ws.Range("B8").Value = 'F003'
$ws.Range("C8").Value = 'Application is installed and launched'
$ws.Range("D8").Value = 'Memory Leak Check - Table Display'
$ws.Range("E8").Value = '1. Navigate to Memory Leak Check section'
$ws.Range("F8").Value = 'Protocol-specific Memory Leak Comparison Table is displayed'
$ws.Range("G8").Value = ''
$ws.Range("H8").Value = ''

# Row 9
$ws.Range("B9").Value = 'F004'
$ws.Range("C9").Value = 'Memory Leak Check section is open'
$ws.Range("D9").Value = 'Memory Leak Check - Comparison'
$ws.Range("E9").Value = '1. Review the Memory Leak Comparison Table<br>2. Compare values to determine memory leak'
$ws.Range("F9").Value = 'User can easily determine if a memory leak has occurred based on the comparison table'
$ws.Range("G9").Value = ''
$ws.Range("H9").Value = ''

# Row 10
$ws.Range("B10").Value = 'F005'
$ws.Range("C10").Value = 'Application is installed and launched'
$ws.Range("D10").Value = 'Debug Log Collection'
$ws.Range("E10").Value = '1. Navigate to Debug Log Collection section<br>2. Click on Run button'
$ws.Range("F10").Value = '1. Script executes successfully<br>2. Logs are collected<br>3. Logs are copied to MFP''s Shared Folder<br>4. Shared Folder opens automatically'
$ws.Range("G10").Value = ''
$ws.Range("H10").Value = ''

# Row 11
$ws.Range("B11").Value = 'F006'
$ws.Range("C11").Value = 'Debug Log Collection failed on first attempt'
$ws.Range("D11").Value = 'Debug Log Collection - Retry'
$ws.Range("E11").Value = '1. Navigate to Debug Log Collection section<br>2. Click on Run button again'
$ws.Range("F11").Value = '1. Script executes successfully<br>2. Logs are collected<br>3. Logs are copied to MFP''s Shared Folder<br>4. Shared Folder opens with logs visible'
$ws.Range("G11").Value = ''
$ws.Range("H11").Value = ''

# Row 12
$ws.Range("B12").Value = 'F007'
$ws.Range("C12").Value = 'Application is installed and launched'
$ws.Range("D12").Value = 'Diagnostic Code Details - ECC'
$ws.Range("E12").Value = '1. Navigate to Diagnostic Code Details section<br>2. Select ECC option'
$ws.Range("F12").Value = 'Relevant job-specific details for ECC are displayed'
$ws.Range("G12").Value = ''
$ws.Range("H12").Value = ''

# Row 13
$ws.Range("B13").Value = 'F008'
$ws.Range("C13").Value = 'Application is installed and launched'
$ws.Range("D13").Value = 'Diagnostic Code Details - Network Protocols'
$ws.Range("E13").Value = '1. Navigate to Diagnostic Code Details section<br>2. Select Network Protocols option'
$ws.Range("F13").Value = 'Relevant job-specific details for Network Protocols are displayed'
$ws.Range("G13").Value = ''
$ws.Range("H13").Value = ''

# Row 14
$ws.Range("B14").Value = 'F009'
$ws.Range("C14").Value = 'Application is installed and launched'
$ws.Range("D14").Value = 'Diagnostic Code Details - High Security Mode'
$ws.Range("E14").Value = '1. Navigate to Diagnostic Code Details section<br>2. Select High Security Mode option'
$ws.Range("F14").Value = 'Relevant job-specific details for High Security Mode are displayed'
$ws.Range("G14").Value = ''
$ws.Range("H14").Value = ''

# Row 15
$ws.Range("B15").Value = 'F010'
$ws.Range("C15").Value = 'Application is installed and launched'
$ws.Range("D15").Value = 'Diagnostic Code Details - Common Codes'
$ws.Range("E15").Value = '1. Navigate to Diagnostic Code Details section<br>2. Select commonly used diagnostic codes'
$ws.Range("F15").Value = 'Relevant job-specific details for selected diagnostic codes are displayed'
$ws.Range("G15").Value = ''
$ws.Range("H15").Value = ''

# Row 16
$ws.Range("B16").Value = 'F011'
$ws.Range("C16").Value = 'Application is installed and launched'
$ws.Range("D16").Value = '08 Diagnostic Code Value - Get'
$ws.Range("E16").Value = '1. Navigate to 08 Diagnostic Code Value section<br>2. Select a diagnostic code<br>3. Click on Get button'
$ws.Range("F16").Value = 'Current value of the selected diagnostic code is displayed'
$ws.Range("G16").Value = ''
$ws.Range("H16").Value = ''

# Row 17
$ws.Range("B17").Value = 'F012'
$ws.Range("C17").Value = 'Application is installed and launched'
$ws.Range("D17").Value = '08 Diagnostic Code Value - Set'
$ws.Range("E17").Value = '1. Navigate to 08 Diagnostic Code Value section<br>2. Select a diagnostic code<br>3. Enter a new value<br>4. Click on Set button'
$ws.Range("F17").Value = 'The diagnostic code value is updated successfully'
$ws.Range("G17").Value = ''
$ws.Range("H17").Value = ''

# Row 18
$ws.Range("B18").Value = 'F013'
$ws.Range("C18").Value = 'Application is installed and launched'
$ws.Range("D18").Value = 'Protocol Configuration - Get'
$ws.Range("E18").Value = '1. Navigate to Protocol Configuration section<br>2. Open Protocol Selection Window<br>3. Select a protocol<br>4. Click on Get button'
$ws.Range("F18").Value = 'Current value of the selected protocol is displayed'
$ws.Range("G18").Value = ''
$ws.Range("H18").Value = ''

# Row 19
$ws.Range("B19").Value = 'F014'
$ws.Range("C19").Value = 'Application is installed and launched'
$ws.Range("D19").Value = 'Protocol Configuration - Set (Pending)'
$ws.Range("E19").Value = '1. Navigate to Protocol Configuration section<br>2. Open Protocol Selection Window<br>3. Select a protocol<br>4. Enter a new value<br>5. Click on Set button'
$ws.Range("F19").Value = 'Feature is marked as "Not Implemented" or similar message is displayed'
$ws.Range("G19").Value = ''
$ws.Range("H19").Value = 'Feature pending implementation as per SRS'

# Row 20
$ws.Range("B20").Value = 'Test Case ID'
$ws.Range("C20").Value = 'Preconditions'
$ws.Range("D20").Value = 'Test Condition'
$ws.Range("E20").Value = 'Steps with description'
$ws.Range("F20").Value = 'Expected Result'
$ws.Range("G20").Value = 'Actual Result'
$ws.Range("H20").Value = 'Remarks'

# Row 21
$ws.Range("B21").Value = '--------------'
$ws.Range("C21").Value = '---------------'
$ws.Range("D21").Value = '----------------'
$ws.Range("E21").Value = '------------------------'
$ws.Range("F21").Value = '-----------------'
$ws.Range("G21").Value = '--------------'
$ws.Range("H21").Value = '---------'

# Row 22
$ws.Range("B22").Value = 'NF001'
$ws.Range("C22").Value = 'Desktop environment with application installed'
$ws.Range("D22").Value = 'Desktop Compatibility'
$ws.Range("E22").Value = '1. Install application using Desktop zip file<br>2. Launch application<br>3. Test all major functions'
$ws.Range("F22").Value = 'Application runs correctly on desktop environment'
$ws.Range("G22").Value = ''
$ws.Range("H22").Value = ''

# Row 23
$ws.Range("B23").Value = 'NF002'
$ws.Range("C23").Value = 'Laptop environment with application installed'
$ws.Range("D23").Value = 'Laptop Compatibility'
$ws.Range("E23").Value = '1. Install application using Laptop zip file<br>2. Launch application<br>3. Test all major functions'
$ws.Range("F23").Value = 'Application runs correctly on laptop environment'
$ws.Range("G23").Value = ''
$ws.Range("H23").Value = ''

# Row 24
$ws.Range("B24").Value = 'NF003'
$ws.Range("C24").Value = 'Application is installed and launched'
$ws.Range("D24").Value = 'Performance - Time Saving'
$ws.Range("E24").Value = '1. Perform a complete diagnostic operation using the tool<br>2. Compare time taken with manual method'
$ws.Range("F24").Value = 'Tool reduces testing time by approximately 80% compared to manual methods'
$ws.Range("G24").Value = ''
$ws.Range("H24").Value = ''

# Row 25
$ws.Range("B25").Value = 'NF004'
$ws.Range("C25").Value = 'Application is installed and launched'
$ws.Range("D25").Value = 'Usability - GUI Intuitiveness'
$ws.Range("E25").Value = '1. Ask a new user to perform basic operations without instructions<br>2. Observe user''s ability to navigate and use features'
$ws.Range("F25").Value = 'User can navigate and use basic features without significant confusion'
$ws.Range("G25").Value = ''
$ws.Range("H25").Value = ''

# Row 26
$ws.Range("B26").Value = 'NF005'
$ws.Range("C26").Value = 'Application is installed and launched'
$ws.Range("D26").Value = 'Usability - Error Handling'
$ws.Range("E26").Value = '1. Deliberately perform incorrect operations<br>2. Observe application''s response'
$ws.Range("F26").Value = 'Application provides clear error messages and doesn''t crash'
$ws.Range("G26").Value = ''
$ws.Range("H26").Value = ''

# Row 27
$ws.Range("B27").Value = 'NF006'
$ws.Range("C27").Value = 'Application is installed and launched'
$ws.Range("D27").Value = 'Performance - Resource Usage'
$ws.Range("E27").Value = '1. Launch application<br>2. Monitor CPU and memory usage during operations'
$ws.Range("F27").Value = 'Application uses reasonable system resources without excessive consumption'
$ws.Range("G27").Value = ''
$ws.Range("H27").Value = ''

# Row 28
$ws.Range("B28").Value = 'NF007'
$ws.Range("C28").Value = 'Application is installed and launched'
$ws.Range("D28").Value = 'Installation Process'
$ws.Range("E28").Value = '1. Extract zip file to preferred location<br>2. Launch application by double-clicking MultiFunctionalToolApplication'
$ws.Range("F28").Value = 'Application installs and launches without errors'
$ws.Range("G28").Value = ''
$ws.Range("H28").Value = ''
